$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------------
# Sheet "preguntas": insert scale_min/scale_max before tipo_informe, insert
# run_id before timestamp, update the (changed) timestamp value, and append
# the new metadata.* / dimension_id columns.
# -------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("preguntas")

# Insert two blank columns for scale_min / scale_max (old Q "tipo_informe" -> S)
$ws1.Columns("Q:R").Insert()
# Insert one blank column for run_id (old T "timestamp", now shifted, -> W)
$ws1.Columns("V:V").Insert()

$lastRow = 23

# Headers for the newly inserted columns
$ws1.Range("Q1").Value = "scale_min"
$ws1.Range("R1").Value = "scale_max"
$ws1.Range("V1").Value = "run_id"

# New trailing headers
$ws1.Range("AA1").Value = "metadata.pipeline_version"
$ws1.Range("AB1").Value = "metadata.pipeline_version_history"
$ws1.Range("AC1").Value = "metadata.run_id"
$ws1.Range("AD1").Value = "metadata.run_id_history"
$ws1.Range("AE1").Value = "metadata.criteria_version"
$ws1.Range("AF1").Value = "metadata.criteria_versions"
$ws1.Range("AG1").Value = "dimension_id"

# Give the new trailing header cells the same look (bold / border / center-top)
# as the rest of row 1 by copy-pasting formats from an existing header cell.
$ws1.Range("P1").Copy()
$ws1.Range("AA1:AG1").PasteSpecial(-4122)

# run_id (new column V) + updated timestamp (shifted column W) for every data row
$ws1.Range("V2:V23").Value = "a1de0ddd622744788ce686c90a8f7dbb"
$ws1.Range("W2:W23").Value = "2025-10-26T22:25:27.187688"

# metadata.pipeline_version / history (constant across all rows)
$ws1.Range("AA2:AA23").Value = "0.1.0"
$ws1.Range("AB2:AB23").Value = '["0.1.0"]'

# metadata.run_id / history (constant across all rows)
$ws1.Range("AC2:AC23").Value = "a1de0ddd622744788ce686c90a8f7dbb"
$ws1.Range("AD2:AD23").Value = '["a1de0ddd622744788ce686c90a8f7dbb"]'

# metadata.criteria_version (constant across all rows) - copy from the
# existing criteria_version column (P) so the text type (not numeric) is
# preserved, then fill in the companion "history" list column.
$ws1.Range("P2:P23").Copy()
$ws1.Range("AE2:AE23").PasteSpecial(-4163)
$ws1.Range("AF2:AF23").Value = '["2024.1"]'

# dimension_id mirrors the dimension_name column (F) for every row
$ws1.Range("F2:F23").Copy()
$ws1.Range("AG2:AG23").PasteSpecial(-4163)

# -------------------------------------------------------------------------
# Sheet "resumen": append raw_score + criteria_version columns.
# -------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("resumen")

$ws2.Range("G1").Value = "raw_score"
$ws2.Range("H1").Value = "criteria_version"

# Match the bold/border/centered header look used by the rest of row 1.
$ws2.Range("F1").Copy()
$ws2.Range("G1:H1").PasteSpecial(-4122)

$ws2.Range("G2:G9").Value = 0

# criteria_version (text) - copy from sheet "preguntas" column P so the
# value keeps its text type instead of being coerced to a number.
$ws1.Range("P2:P9").Copy()
$ws2.Range("H2:H9").PasteSpecial(-4163)

# -------------------------------------------------------------------------
# Sheet "indice_global": move criteria_version earlier (right after
# normalized_max), refresh run_id/timestamp/extra_config, and append the
# new extra_criteria_hash column.
# -------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("indice_global")

# Make room for criteria_version right after normalized_max (column F)
$ws3.Columns("G:G").Insert()

# Move the (now shifted to Q) criteria_version column into the new G slot
$ws3.Range("Q1:Q2").Copy()
$ws3.Range("G1:G2").PasteSpecial(-4163)
$ws3.Columns("Q:Q").Delete()

# Refresh run_id / timestamp / extra_config with their new values
$ws3.Range("M2").Value = "a1de0ddd622744788ce686c90a8f7dbb"
$ws3.Range("Q2").Value = "2025-10-26T22:25:27.187688"
$ws3.Range("R2").Value = "{'model_name': 'gpt-4o-mini', 'retries': 2, 'backoff_factor': 2.0, 'timeout_seconds': 60.0, 'prompt_batch_size': 1, 'log_level': 'INFO', 'log_file': None, 'ai_provider': 'mock', 'run_id': 'a1de0ddd622744788ce686c90a8f7dbb', 'document_id': None, 'extra_instructions': None, 'splitter_log_level': 'info', 'splitter_normalize_newlines': True}"

# New trailing column: extra_criteria_hash
$ws3.Range("V1").Value = "extra_criteria_hash"
$ws3.Range("U1").Copy()
$ws3.Range("V1").PasteSpecial(-4122)
$ws3.Range("V2").Value = "5c7fdfed71a0c2f52a7171577c792d9c236420ba3b160ba8e2305c65aa06bc0c"
